# "THE BD IS REAL"
#
# The sheet used to have a header row (idQuestion / idReponse, driven by
# shared strings) followed by 136 data rows (34 questions x 4 reponses).
# The edit removes that header row entirely: every data row shifts up by
# one, the sheet now starts directly with numeric data on row 1, the
# shared-strings table becomes empty (nothing references strings any
# more), and the old "centered" data format (previously applied to rows
# 2:137) now naturally becomes the format of rows 1:136.
#
# Deleting the whole row 1 (rather than rewriting cell-by-cell) is what
# reproduces this precisely: Excel shifts rows 2..137 up into 1..136,
# keeps their existing per-cell formatting, recomputes the used range to
# A1:B136, and drops the now-unused "idQuestion"/"idReponse" shared
# strings automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

[void]$ws.Rows.Item(1).Delete()

# Mirror the author's final selection: the whole first row is selected
# (sqref A1:XFD1) rather than a single cell.
[void]$ws.Rows.Item(1).Select()
